$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = 7.8
$ws.Range("D4").Value = 8.2
$ws.Range("D5").Value = 9.4
$ws.Range("D7").Value = 9.3
$ws.Range("D9").Value = 7.7
$ws.Range("D11").Value = 6.2
$ws.Range("D12").Value = 5.2
$ws.Range("D15").Value = 2.5
$ws.Range("D20").Value = 7.8
$ws.Range("D21").Value = 8.1
$ws.Range("D22").Value = 9.2
$ws.Range("D23").Value = 8.9
$ws.Range("D24").Value = 8.5
$ws.Range("D25").Value = 7.9
$ws.Range("D26").Value = 8.3
$ws.Range("D30").Value = 4.1
$ws.Range("D31").Value = 3.4
$ws.Range("D36").Value = 7.1
$ws.Range("D37").Value = 8.2
$ws.Range("D38").Value = 7.2
$ws.Range("D39").Value = 7.9
$ws.Range("D40").Value = 8.5
$ws.Range("D41").Value = 7.8
$ws.Range("D42").Value = 8
$ws.Range("D43").Value = 7.3
$ws.Range("D44").Value = 7.5
$ws.Range("D45").Value = 6.9
$ws.Range("D46").Value = 5.5
$ws.Range("D47").Value = 5.4
$ws.Range("D48").Value = 4.3
$ws.Range("D52").Value = 2
